$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell D4:E4 ("Sinais"), mirroring the A4:B4 header format ---
#     Merge the still-blank range first (before it holds any content) so the
#     format-only paste below lands cleanly on both cells without Excel
#     splitting the border into "outer edge only" variants.
$ws.Range("D4:E4").Merge()
$ws.Range("A4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Duplicate the "numbering" column (A5:A16, centered fill) down column D,
#     now extended to 16 rows (D5:D20) -------------------------------------
$ws.Range("A5:A16").Copy()
$ws.Range("D5:D16").PasteSpecial(-4122)
$ws.Range("A5:A8").Copy()
$ws.Range("D17:D20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Duplicate the "student id" column (B5:B16, alternating fill) down
#     column E, also extended to 16 rows (E5:E20) ---------------------------
$ws.Range("B5:B16").Copy()
$ws.Range("E5:E16").PasteSpecial(-4122)
$ws.Range("B5:B8").Copy()
$ws.Range("E17:E20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the actual values --------------------------------------------
$ws.Range("D4").Value = "Sinais"

for ($i = 0; $i -lt 16; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 4).Value = $i + 1
    if (($i % 2) -eq 0) {
        $ws.Cells.Item($row, 5).Value = 1211155
    } else {
        $ws.Cells.Item($row, 5).Value = 1210957
    }
}

# --- Update the active selection -------------------------------------------
$ws.Range("G19").Select()
